# Actualización de minutas de la reunion 1
# Slide 8 ("Minutas y acciones (independiente de las tareas)"),
# Content Placeholder shape - "Campos:" bullet list.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

$GREEN = 65280   # RGB(0,255,0)   -> srgbClr 00FF00
$YELLOW = 65535  # RGB(255,255,0) -> srgbClr FFFF00

# --- Paragraph 5: "Pasar info y repo pagina a Bryan " ---------------------
# Update the trailing run's text, then highlight the whole paragraph green.
$para5 = $tr.Paragraphs(5)
$bryanTailStart = $para5.Start + 10   # after "Pasar " (6) + "info" (4)
$bryanTailLen = $para5.Length - 10
$bryanTail = $tr.Characters($bryanTailStart, $bryanTailLen)
$bryanTail.Text = " y repo pagina a Bryan – va a comenzar una nueva"

# re-fetch paragraph (its length changed after the text edit)
$para5 = $tr.Paragraphs(5)
$para5.Font.Highlight.RGB = $GREEN

# --- Paragraph 6: "Pasar información de MCU a Luis " ----------------------
# Already highlighted green; keep it green (no visible change expected).
$para6 = $tr.Paragraphs(6)
$para6.Font.Highlight.RGB = $GREEN

# --- Paragraph 7: "Ver con Gabo el diseño ... sensado" --------------------
# Highlight yellow.
$para7 = $tr.Paragraphs(7)
$para7.Font.Highlight.RGB = $YELLOW

# --- Paragraph 8: "Re organizar Excel de gastos para pasar a Quiroz" ------
# Highlight green.
$para8 = $tr.Paragraphs(8)
$para8.Font.Highlight.RGB = $GREEN
